# Apply the edits described by the commit diff:
#  - B6: fix typo "I diagree with the news story. " -> "I disagree with the news story. "
#  - B8: "I also agree. " -> "I also agree with the news story. "
#  - B14: drop the appended duplicate sentence, leaving just
#         "I can fix it when a thermostat is not working correctly."
#  - Update the sheet view: scroll so column B is the left-most visible
#    column, and change the active selection from B24 to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value  = "I disagree with the news story. "
$ws.Range("B8").Value  = "I also agree with the news story. "
$ws.Range("B14").Value = "I can fix it when a thermostat is not working correctly."

# Update the view: scroll the window so column B is left-most visible,
# then move the selection to B8.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("B8").Select()
